# Updates FY Future Years to align with IT for EU model
# The "About" sheet's Initial Time (A8) drives all the Future-Year labels
# on the "FY" sheet via formulas. Bumping it from 2021 to 2022 cascades
# the recalculation through the whole FY column.

$wb = $excel.ActiveWorkbook
$about = $wb.Worksheets.Item("About")

$about.Range("A8").Value = 2022
